# Started reviewing WOS returns:
# Insert three new review columns (Compare With/Without Climate?, OOS Score
# with Climate?, Note) in front of the existing WOS export columns, and
# begin filling in the review for the first ~18 records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert three blank columns at the front; this pushes the existing
#     Authors/Title/Journal/Year/DOI columns from A:E to D:H ---
$ws.Range("A1:C1").EntireColumn.Insert()

# --- New header row ---
$ws.Range("A1").Value = "Compare With and Without Climate?"
$ws.Range("B1").Value = "OOS Score with Climate?"
$ws.Range("C1").Value = "Note"

# --- Fill in the review columns for the records looked at so far ---
$ws.Range("A2").Value = "N"
$ws.Range("B2").Value = "N"

$ws.Range("A3").Value = "N"
$ws.Range("B3").Value = "N"

$ws.Range("A4").Value = "N"
$ws.Range("B4").Value = "N"

$ws.Range("A5").Value = "N"
$ws.Range("B5").Value = "Y"

$ws.Range("A6").Value = "Y"
$ws.Range("B6").Value = "Y"

$ws.Range("A7").Value = "N"
$ws.Range("B7").Value = "N"

$ws.Range("A8").Value = "N"
$ws.Range("B8").Value = "Y"

$ws.Range("A9").Value = "N"
$ws.Range("B9").Value = "N"

$ws.Range("A10").Value = "N"
$ws.Range("B10").Value = "N"

$ws.Range("A11").Value = "N"
$ws.Range("B11").Value = "Y"

$ws.Range("A12").Value = "N"
$ws.Range("B12").Value = "N"

$ws.Range("A13").Value = "N"
$ws.Range("B13").Value = "N"

$ws.Range("A14").Value = "Y"
$ws.Range("B14").Value = "Y"

$ws.Range("A15").Value = "N"
$ws.Range("B15").Value = "Y"

$ws.Range("A16").Value = "Y"
$ws.Range("B16").Value = "Y"
$ws.Range("C16").Value = "not out of sample, AUC only"

$ws.Range("A17").Value = "N"
$ws.Range("B17").Value = "Y"
$ws.Range("C17").Value = "not oos"

$ws.Range("A18").Value = "N"
$ws.Range("B18").Value = "N"

$ws.Range("A19").Value = "N"
$ws.Range("B19").Value = "N"

# --- Column widths for the new columns (A-D); E-H keep their original
#     widths automatically since they were shifted with the data ---
$ws.Columns.Item(1).ColumnWidth = 41.1666667
$ws.Columns.Item(2).ColumnWidth = 27.6666667
$ws.Columns.Item(3).ColumnWidth = 27.6666667
$ws.Columns.Item(4).ColumnWidth = 13.1666667

# --- Defined name now needs to cover the shifted data columns ---
$wb.Names.Item("savedrecs").RefersTo = "=Sheet1!`$D`$2:`$I`$120"

# --- Update the active selection to where review work left off ---
$ws.Range("E20").Select()
